$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix total marks error: update "Marking" row (B11 / C11) and "Total" row (B12 / E12)
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

$ws.Range("B12").Value = 68
$ws.Range("E12").Value = "68 / 112"
